$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 461
$ws.Range("I33").Value = 425.06668
$ws.Range("K33").Value = 425.06668
$ws.Range("M33").Value = -196.06668
$ws.Range("H70").Value = 16767967
$ws.Range("I70").Value = 83833980
$ws.Range("J70").Value = 1462.5
$ws.Range("K70").Value = 251501940
$ws.Range("L70").Value = 4387.5
$ws.Range("M70").Value = -251501670
$ws.Range("N70").Value = -4927.5
$ws.Range("H73").Value = 16767967
$ws.Range("I73").Value = 83833980
$ws.Range("J73").Value = 1462.5
$ws.Range("K73").Value = 251501940
$ws.Range("L73").Value = 4387.5
$ws.Range("M73").Value = -251501004
$ws.Range("N73").Value = -6259.5
$ws.Range("H112").Value = 1671.1818
$ws.Range("J112").Value = 1745.9667
$ws.Range("L112").Value = 5237.9001
$ws.Range("N112").Value = -7453.9001
$ws.Range("H132").Value = 3667.8262
$ws.Range("I132").Value = 3775.6428
$ws.Range("K132").Value = 11326.9284
$ws.Range("M132").Value = -8796.928400000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 0
$ws.Range("K6").Value = 0
$ws.Range("M6").ClearContents()
$ws.Range("H45").Value = 1173.6
$ws.Range("I45").Value = 1132.2858
$ws.Range("J45").Value = 1270
$ws.Range("K45").Value = 1132.2858
$ws.Range("L45").Value = 1270
$ws.Range("M45").Value = -755.2858000000001
$ws.Range("N45").Value = -2024
$ws.Range("H119").Value = 23270.688
$ws.Range("J119").Value = 23270.688
$ws.Range("L119").Value = 23270.688
$ws.Range("N119").Value = -32946.68799999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 201.42857
$ws.Range("I19").Value = 201.42857
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 201.42857
$ws.Range("L19").Value = 0
$ws.Range("M19").Value = -31.42857000000001
$ws.Range("N19").ClearContents()
$ws.Range("H24").Value = 201.42857
$ws.Range("I24").Value = 201.42857
$ws.Range("J24").Value = 0
$ws.Range("K24").Value = 201.42857
$ws.Range("L24").Value = 0
$ws.Range("M24").Value = -31.42857000000001
$ws.Range("N24").ClearContents()
$ws.Range("H99").Value = 2232.5833
$ws.Range("I99").Value = 2128
$ws.Range("J99").Value = 2379
$ws.Range("K99").Value = 2128
$ws.Range("L99").Value = 2379
$ws.Range("M99").Value = -630
$ws.Range("N99").Value = -5375
$ws.Range("H122").Value = 2396.6
$ws.Range("I122").Value = 2116
$ws.Range("J122").Value = 3051.3333
$ws.Range("K122").Value = 6348
$ws.Range("L122").Value = 9153.999899999999
$ws.Range("M122").Value = -3898
$ws.Range("N122").Value = -14053.9999
$ws.Range("H126").Value = 2232.5833
$ws.Range("I126").Value = 2128
$ws.Range("J126").Value = 2379
$ws.Range("K126").Value = 6384
$ws.Range("L126").Value = 7137
$ws.Range("M126").Value = -3914
$ws.Range("N126").Value = -12077
$ws.Range("H134").Value = 3659.5833
$ws.Range("I134").Value = 2525.5
$ws.Range("J134").Value = 4226.625
$ws.Range("K134").Value = 7576.5
$ws.Range("L134").Value = 12679.875
$ws.Range("M134").Value = -5041.5
$ws.Range("N134").Value = -17749.875

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 1017084.7
$ws.Range("J12").Value = 1288278.1
$ws.Range("L12").Value = 3864834.3
$ws.Range("N12").Value = -3865180.3
$ws.Range("H113").Value = 664.7308
$ws.Range("I113").Value = 599.8823
$ws.Range("J113").Value = 787.2222
$ws.Range("K113").Value = 1799.6469
$ws.Range("L113").Value = 2361.6666
$ws.Range("M113").Value = 370.3531
$ws.Range("N113").Value = -6701.6666
$ws.Range("H120").Value = 15293.066
$ws.Range("I120").Value = 14000
$ws.Range("J120").Value = 15385.429
$ws.Range("K120").Value = 42000
$ws.Range("L120").Value = 46156.287
$ws.Range("M120").Value = -37162
$ws.Range("N120").Value = -55832.287
$ws.Range("H122").Value = 348.33334
$ws.Range("I122").Value = 352.72726
$ws.Range("J122").Value = 300
$ws.Range("K122").Value = 3174.54534
$ws.Range("L122").Value = 2700
$ws.Range("M122").Value = -724.5453400000001
$ws.Range("N122").Value = -7600
$ws.Range("H123").Value = 1988.3334
$ws.Range("I123").Value = 1988.3334
$ws.Range("J123").Value = 0
$ws.Range("K123").Value = 5965.0002
$ws.Range("L123").Value = 0
$ws.Range("M123").Value = -3515.0002
$ws.Range("N123").ClearContents()
$ws.Range("H125").Value = 3914.4443
$ws.Range("J125").Value = 3914.4443
$ws.Range("L125").Value = 11743.3329
$ws.Range("N125").Value = -21583.3329
$ws.Range("H133").Value = 3675.818
$ws.Range("I133").Value = 1663.8462
$ws.Range("K133").Value = 4991.5386
$ws.Range("M133").Value = 68.46140000000014

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H62").Value = 23000
$ws.Range("I62").Value = 23000
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 23000
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -22314
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 23000
$ws.Range("I65").Value = 23000
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 69000
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -65568
$ws.Range("N65").ClearContents()
$ws.Range("H113").Value = 1437.8235
$ws.Range("I113").Value = 894.3333
$ws.Range("J113").Value = 2049.25
$ws.Range("K113").Value = 894.3333
$ws.Range("L113").Value = 2049.25
$ws.Range("M113").Value = 1275.6667
$ws.Range("N113").Value = -6389.25
$ws.Range("H122").Value = 2449.5334
$ws.Range("I122").Value = 1991.8182
$ws.Range("K122").Value = 5975.4546
$ws.Range("M122").Value = -3525.4546

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H25").Value = 334999.66
$ws.Range("I25").Value = 334999.66
$ws.Range("K25").Value = 334999.66
$ws.Range("M25").Value = -334769.66
$ws.Range("H61").Value = 2960
$ws.Range("I61").Value = 1650
$ws.Range("J61").Value = 3833.3333
$ws.Range("K61").Value = 1650
$ws.Range("L61").Value = 3833.3333
$ws.Range("M61").Value = -1448
$ws.Range("N61").Value = -4237.3333
$ws.Range("H69").Value = 193354.33
$ws.Range("J69").Value = 193354.33
$ws.Range("L69").Value = 193354.33
$ws.Range("N69").Value = -194976.33
$ws.Range("H72").Value = 193354.33
$ws.Range("J72").Value = 193354.33
$ws.Range("L72").Value = 580062.99
$ws.Range("N72").Value = -588174.99
$ws.Range("H113").Value = 2960
$ws.Range("I113").Value = 1650
$ws.Range("J113").Value = 3833.3333
$ws.Range("K113").Value = 1650
$ws.Range("L113").Value = 3833.3333
$ws.Range("M113").Value = 520
$ws.Range("N113").Value = -8173.3333
$ws.Range("H119").Value = 45000
$ws.Range("J119").Value = 45000
$ws.Range("L119").Value = 45000
$ws.Range("N119").Value = -54676
$ws.Range("H136").Value = 25252148
$ws.Range("I136").Value = 34484212
$ws.Range("J136").Value = 913073.6
$ws.Range("K136").Value = 103452636
$ws.Range("L136").Value = 2739220.8
$ws.Range("M136").Value = -103450086
$ws.Range("N136").Value = -2744320.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1000
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 1000
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 3000
$ws.Range("M107").ClearContents()
$ws.Range("N107").Value = -6840
$ws.Range("H136").Value = 2287.875
$ws.Range("I136").Value = 2160.6
$ws.Range("J136").Value = 2500
$ws.Range("K136").Value = 6481.799999999999
$ws.Range("L136").Value = 7500
$ws.Range("M136").Value = -3931.799999999999
$ws.Range("N136").Value = -12600
